# Generate Report for Handback
# Refresh the handback-status timestamps to reflect the latest report run.

$wb = $excel.ActiveWorkbook

# "Overview" sheet — Latest HO Xliff Generate Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value2 = "2016-09-07 17:35:56"

# "zh-cn" sheet — Correspond Handoff / Handback Datetime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value2 = "2016-09-07 17:35:51"
$wsZhCn.Range("K2").Value2 = "2016-09-07 17:36:33"

# "de-de" sheet — Correspond Handoff / Handback Datetime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value2 = "2016-09-07 17:35:56"
$wsDeDe.Range("K2").Value2 = "2016-09-07 17:36:42"
